# Update scripts with new TPM values (Nell2-Robo3, FAPs -> ECs / FAPs->FAPs rows)
# and drop the obsolete FAPs -> MuSCs row entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs): refreshed receptor/edge statistics.
$ws.Range("M2").Value = 3.555949666666667
$ws.Range("N2").Value = 10.667849
$ws.Range("O2").Value = 0.9961127894270114
$ws.Range("P2").Value = 0.9961127894270114
$ws.Range("Q2").Value = 1.347633804673334
$ws.Range("R2").Value = 12.12870424206
$ws.Range("S2").Value = 0.9961127894270114
$ws.Range("T2").Value = 0.9961127894270114

# Row 3 (FAPs -> FAPs): refreshed derived-specificity columns.
$ws.Range("O3").Value = 0.003887210572988658
$ws.Range("P3").Value = 0.003887210572988658
$ws.Range("S3").Value = 0.003887210572988658
$ws.Range("T3").Value = 0.003887210572988658

# Row 4 (FAPs -> MuSCs) no longer exists in the refreshed output; remove it
# so row 3 becomes the last data row and "MuSCs" drops out of the shared
# string table along with it.
$ws.Rows(4).Delete()
